# Disable investment in Storage (in Excel)
# Column R ("EnableInvest") for the Storage rows (8-12) is set from 1 to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

$ws.Range("R8:R12").Value = 0

# Leave the selection where it ended up after editing the last cell (R12 -> R13)
$ws.Range("R13").Select() | Out-Null
